# Add data for 2022-12-20 (update "through December 11" -> "through December 12"
# column, updating the header label, the sheet/tab name, and the carjacking
# counts for that column across the affected neighborhoods).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab and update the header cell that describes the date range.
$ws.Name = "Through 2022-12-12"
$ws.Range("B1").Value = "December 2022 (through December 12)"

# Cell updates: Row, Column (1-based), New Value
$updates = @(
    @(2,   2,  1),
    @(4,  74,  1),
    @(7,  26,  4),
    @(7,  74,  1),
    @(8,  14,  1),
    @(9,  26,  2),
    @(10, 26,  1),
    @(14, 50,  4),
    @(18, 74,  2),
    @(19, 14,  3),
    @(19, 74,  1),
    @(20, 26,  8),
    @(20, 86,  1),
    @(26, 14,  2),
    @(27, 38,  2),
    @(28, 50,  2),
    @(28, 62,  5),
    @(36, 86,  3),
    @(40,  2,  1),
    @(41, 14,  2),
    @(41, 38,  1),
    @(48, 14,  1),
    @(52,  2,  1),
    @(57, 14,  3),
    @(92, 50,  1),
    @(96, 26,  1)
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $ws.Cells.Item($row, $col).Value = $val
}
